# Auto-generated edit script: updates F/G column numeric values
# per the commit diff (gh-pages output regenerated at 456a3b4).
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1725
$ws.Range("F3").Value = 10058
$ws.Range("F6").Value = 590
$ws.Range("F8").Value = 1586
$ws.Range("F10").Value = 364
$ws.Range("F12").Value = 194
$ws.Range("F16").Value = 125
$ws.Range("F18").Value = 10
$ws.Range("F19").Value = 82
$ws.Range("F23").Value = 95
$ws.Range("F24").Value = 869
$ws.Range("F25").Value = 684
$ws.Range("F27").Value = 33
$ws.Range("F28").Value = 433
$ws.Range("F29").Value = 219
$ws.Range("F30").Value = 98
$ws.Range("F31").Value = 343
$ws.Range("F33").Value = 365
$ws.Range("F34").Value = 520
$ws.Range("F35").Value = 558
$ws.Range("F36").Value = 719
$ws.Range("F39").Value = 805
$ws.Range("F40").Value = 370
$ws.Range("F41").Value = 325
$ws.Range("F42").Value = 3
$ws.Range("F43").Value = 349
$ws.Range("F45").Value = 347
$ws.Range("F46").Value = 76

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 44
$ws.Range("F8").Value = 70
$ws.Range("F13").Value = 97
$ws.Range("G16").Value = 180
$ws.Range("F19").Value = 555
$ws.Range("F20").Value = 1094
$ws.Range("F24").Value = 8
$ws.Range("F27").Value = 361
$ws.Range("F33").Value = 154
$ws.Range("F35").Value = 42
$ws.Range("F38").Value = 129
$ws.Range("F41").Value = 35

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 801
$ws.Range("F6").Value = 2495
$ws.Range("F7").Value = 4033
$ws.Range("F8").Value = 52
$ws.Range("F10").Value = 271
$ws.Range("F11").Value = 177

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1725
$ws.Range("F3").Value = 801
$ws.Range("F4").Value = 10058
$ws.Range("F7").Value = 4033
$ws.Range("F8").Value = 52
$ws.Range("F9").Value = 271
$ws.Range("F10").Value = 271
$ws.Range("F11").Value = 590
$ws.Range("F12").Value = 1586
$ws.Range("F14").Value = 364
$ws.Range("F15").Value = 194
$ws.Range("F19").Value = 125
$ws.Range("F22").Value = 97
$ws.Range("F24").Value = 82
$ws.Range("F28").Value = 1094
$ws.Range("F30").Value = 872
$ws.Range("F31").Value = 684
$ws.Range("F34").Value = 361
$ws.Range("F35").Value = 343
$ws.Range("F37").Value = 365
$ws.Range("F38").Value = 520
$ws.Range("F39").Value = 558
$ws.Range("F41").Value = 719
$ws.Range("F43").Value = 805
$ws.Range("F44").Value = 370
$ws.Range("F45").Value = 42
$ws.Range("F46").Value = 325
$ws.Range("F47").Value = 129
$ws.Range("F48").Value = 349
$ws.Range("F49").Value = 347
